$wb = $excel.ActiveWorkbook

# --- Login sheet: add "pageheader" / "Swag Labs" column C ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("C1").Value = "pageheader"
for ($r = 2; $r -le 7; $r++) {
    $wsLogin.Cells.Item($r, 3).Value = "Swag Labs"
}
$wsLogin.Range("A1:C2").Select() | Out-Null

# --- Rename Sheet2 -> FilterProduct and populate it ---
$wsFilter = $wb.Worksheets.Item("Sheet2")
$wsFilter.Name = "FilterProduct"

$wsFilter.Range("A1").Value = "username"
$wsFilter.Range("B1").Value = "password"
$wsFilter.Range("C1").Value = "pageheader"
$wsFilter.Range("D1").Value = "filterproduct"

$wsFilter.Range("A2").Value = "standard_user"
$wsFilter.Range("B2").Value = "secret_sauce"
$wsFilter.Range("C2").Value = "Swag Labs"
$wsFilter.Range("D2").Value = "Price (low to high)"

$wsFilter.Range("D3").Select() | Out-Null
$wsFilter.Activate() | Out-Null

Write-Output "done"
